# Update the "想去人数" (want-to-go count) figures that changed between
# consecutive generated-data snapshots of this ACG-convention tracker.
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both list the
# same set of exhibition events (just interleaved with other event types in
# the latter), so the same seven F-column values are bumped in both sheets.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F14").Value = 4351
$wsExhibitions.Range("F16").Value = 1679
$wsExhibitions.Range("F22").Value = 995
$wsExhibitions.Range("F23").Value = 298
$wsExhibitions.Range("F35").Value = 1694
$wsExhibitions.Range("F40").Value = 616
$wsExhibitions.Range("F41").Value = 313

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F16").Value = 4351
$wsAllTypes.Range("F18").Value = 1679
$wsAllTypes.Range("F26").Value = 995
$wsAllTypes.Range("F27").Value = 298
$wsAllTypes.Range("F35").Value = 1694
$wsAllTypes.Range("F42").Value = 616
$wsAllTypes.Range("F43").Value = 313
